$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.480.54'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '1.914.48'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.73'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4836'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2896'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.60%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06722'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '109.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.72%  '

$ws.Range("D12").Value = '1.913.94'
$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("E13").Value = '  -0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.276'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6727'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.51%  '

$ws.Range("E16").Value = '  -3.94%  '

$ws.Range("D17").Value = '30.496.50'
$ws.Range("E17").Value = '  +0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007567'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.18%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.508'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.60%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.165.62'
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.480'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.449'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.63%  '

$ws.Range("E29").Value = '  -1.52%  '

$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.158'
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = '  +2.70%  '

$ws.Range("E33").Value = '  -0.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7307'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.98%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.18%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.728'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02030'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.662'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '110.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.017'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4443'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8651'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.812'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9996'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.356'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.280'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.88%  '

$ws.Range("E51").Value = '  +6.27%  '
